$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.610.47'
$ws.Range("E2").Value = '  -0.30%  '
$ws.Range("D3").Value = '1.844.59'
$ws.Range("E3").Value = '  -0.33%  '
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.56'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.89%  '
$ws.Range("E6").Value = '  +0.25%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4308'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.89%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3686'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +1.45%  '
$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07326'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.24%  '
$ws.Range("B10").Value = 'Polygon'
$ws.Range("C10").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8742'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.27%  '
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.98'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +1.58%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.860.03'
$ws.Range("E12").Value = '  -0.57%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.472'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.65%  '
$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.594'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +1.09%  '
$ws.Range("B15").Value = 'TRON'
$ws.Range("C15").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06942'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +0.75%  '
$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.008'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.43%  '
$ws.Range("B17").Value = 'Litecoin'
$ws.Range("C17").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '81.40'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +1.92%  '
$ws.Range("B18").Value = 'ShibaInu'
$ws.Range("C18").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009070'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.44%  '
$ws.Range("B19").Value = 'Dai'
$ws.Range("C19").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.003'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.56'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +1.30%  '
$ws.Range("B21").Value = 'WrappedBTC'
$ws.Range("C21").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D21").Value = '27.718.17'
$ws.Range("E21").Value = '  +0.03%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.085'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +2.31%  '
$ws.Range("B23").Value = 'Cosmos'
$ws.Range("C23").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.98'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +5.69%  '
$ws.Range("B24").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C24").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D24").Value = '2.095.41'
$ws.Range("E24").Value = '  -1.08%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.988'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.20%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.35'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +0.20%  '
$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.00'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.99%  '
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.328'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.94%  '
$ws.Range("B29").Value = 'BitcoinCash'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '115.94'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -4.55%  '
$ws.Range("B30").Value = 'LidoDAOToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.879'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.79%  '
$ws.Range("B31").Value = 'Stellar'
$ws.Range("C31").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08903'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.24%  '
$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7841'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +3.07%  '
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.607'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +1.96%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.986'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +0.57%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.165'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +5.65%  '
$ws.Range("B36").Value = 'Frax'
$ws.Range("C36").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.004'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +0.31%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.109'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.68%  '
$ws.Range("B38").Value = 'Hedera'
$ws.Range("C38").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05435'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +0.87%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01963'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.62%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.828'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.21%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5174'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +1.83%  '
$ws.Range("B42").Value = 'Algorand'
$ws.Range("C42").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1693'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +2.32%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.761'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.13%  '
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.630'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +3.29%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.67'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +3.45%  '
$ws.Range("B46").Value = 'PaxosStandard'
$ws.Range("C46").Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.006'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -32.65%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4789'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.31%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '106.55'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +1.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06543'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.004'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.30%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.663'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.68%  '
